$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 99 (existing rows 99:111 shift down to 100:112)
$ws.Rows("99:99").Insert()

# Populate the newly inserted row 99 with the new weekly price record
$ws.Range("A99").Value = 7
$ws.Range("B99").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C99").Value = "Ñuble"
$ws.Range("D99").Value = 45154
$ws.Range("E99").Value = 16
$ws.Range("F99").Value = 100112044
$ws.Range("G99").Value = "Perejil"
$ws.Range("H99").Value = "Sin especificar"
$ws.Range("I99").Value = "Primera"
$ws.Range("J99").Value = 200
$ws.Range("K99").Value = 1500
$ws.Range("L99").Value = 1500
$ws.Range("M99").Value = 1500
$ws.Range("N99").Value = "`$/atado 0,5 a 1 kilo"
$ws.Range("O99").Value = "Región de Ñuble"
$ws.Range("P99").Value = 1500
$ws.Range("Q99").Value = 1
$ws.Range("R99").Value = "Hortaliza"
